# Apply crypto price/volume updates as produced by the GitHub Actions
# "Updated cryptos list" automation.
#
# Target cells were originally stored as inline-string text (t="inlineStr")
# with no custom cell style. Some of the new values look like numbers
# (e.g. "1.000", "18.00", "0.9990") and Excel would normally auto-convert
# a plain .Value assignment of such a string into a numeric cell, losing
# the exact textual representation (trailing zeros, etc.) and changing the
# cell type. To avoid that, we briefly force the cell's number format to
# Text ("@") before assigning the value, then restore the cell style back
# to "Normal" afterwards so the resulting cell keeps a plain/default style
# (matching the source workbook) while still holding the exact text value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @('D2', '29.090.25'),
    @('D3', '1.843.22'),
    @('E3', '  -2.11%  '),
    @('D4', '1.000'),
    @('E4', '  +0.08%  '),
    @('D5', '0.6995'),
    @('E5', '  -5.80%  '),
    @('D6', '236.94'),
    @('E6', '  -2.38%  '),
    @('E7', '  -0.14%  '),
    @('D8', '0.3026'),
    @('E8', '  -4.24%  '),
    @('D9', '0.07383'),
    @('E9', '  +2.24%  '),
    @('D10', '23.28'),
    @('E10', '  -6.31%  '),
    @('D11', '0.08106'),
    @('E11', '  -2.94%  '),
    @('B12', 'Polygon'),
    @('C12', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'),
    @('D12', '0.7230'),
    @('E12', '  -4.29%  '),
    @('B13', 'WrappedEther'),
    @('C13', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'),
    @('D13', '1.822.73'),
    @('E13', '  -5.09%  '),
    @('D14', '5.199'),
    @('E14', '  -3.78%  '),
    @('D15', '88.88'),
    @('E15', '  -3.84%  '),
    @('D16', '29.246.06'),
    @('E16', '  -2.68%  '),
    @('D17', '5.767'),
    @('E17', '  -6.21%  '),
    @('D18', '241.03'),
    @('E18', '  -3.21%  '),
    @('D19', '0.000007646'),
    @('E19', '  -2.87%  '),
    @('D20', '12.97'),
    @('E20', '  -4.64%  '),
    @('D21', '1.000'),
    @('E21', '  -0.07%  '),
    @('D22', '2.117.64'),
    @('E22', '  -1.24%  '),
    @('D23', '0.9990'),
    @('E23', '  -0.16%  '),
    @('D24', '7.574'),
    @('E24', '  -5.77%  '),
    @('D25', '0.1470'),
    @('E25', '  -5.96%  '),
    @('D26', '161.87'),
    @('E26', '  -2.60%  '),
    @('D27', '8.937'),
    @('E27', '  -4.10%  '),
    @('D28', '18.00'),
    @('E28', '  -3.78%  '),
    @('D29', '1.930'),
    @('E29', '  -5.39%  '),
    @('D30', '1.374'),
    @('E30', '  -8.58%  '),
    @('D31', '4.441'),
    @('E31', '  -3.51%  '),
    @('E32', '  -3.27%  '),
    @('D33', '3.996'),
    @('E33', '  -5.37%  '),
    @('D34', '0.05194'),
    @('E34', '  -3.33%  '),
    @('D35', '1.180'),
    @('E35', '  -5.80%  '),
    @('D36', '0.7091'),
    @('E36', '  -6.20%  '),
    @('D37', '1.002'),
    @('E37', '  -1.02%  '),
    @('D38', '2.646'),
    @('E38', '  -2.24%  '),
    @('E39', '  -5.23%  '),
    @('D40', '2.669'),
    @('E40', '  -3.28%  '),
    @('D41', '0.9027'),
    @('E41', '  +4.95%  '),
    @('D42', '0.4278'),
    @('E42', '  -6.05%  '),
    @('D43', '5.889'),
    @('E43', '  -4.51%  '),
    @('B44', 'Maker'),
    @('C44', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'),
    @('D44', '1.048.78'),
    @('E44', '  -5.72%  '),
    @('B45', 'Aave'),
    @('C45', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'),
    @('D45', '69.82'),
    @('E45', '  -4.28%  '),
    @('D46', '0.9992'),
    @('E46', '  -0.09%  '),
    @('D47', '101.45'),
    @('E47', '  -3.37%  '),
    @('D48', '1.749'),
    @('E48', '  -6.76%  '),
    @('D49', '7.095'),
    @('E49', '  -7.00%  '),
    @('D50', '1.983.88'),
    @('E50', '  -3.65%  '),
    @('D51', '9.179'),
    @('E51', '  -4.04%  ')
)

foreach ($entry in $updates) {
    $ref = $entry[0]
    $val = $entry[1]
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}
